# Fruta / hortaliza, semanal
# Updates weekly price data: columns D (Fecha, date serial), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) get reshuffled across rows 2-9 (row 7 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (taken from the diff target state)
$data = @{
    2 = @{ D = 44589; J = 110; K = 5000; L = 6000; M = 5500; P = 92 }
    3 = @{ D = 44382; J = 160; K = 7000; L = 8000; M = 7438; P = 124 }
    4 = @{ D = 44242; J = 160; K = 5000; L = 5500; M = 5250; P = 88 }
    5 = @{ D = 44494; J = 120; K = 5000; L = 6000; M = 5500; P = 92 }
    6 = @{ D = 44281; J = 120; K = 5500; L = 6000; M = 5750; P = 96 }
    7 = @{ D = 44400; J = 120; K = 9000; L = 10000; M = 9500; P = 158 }
    8 = @{ D = 44362; J = 120; K = 8000; L = 9000; M = 8500; P = 142 }
    9 = @{ D = 44421; J = 100; K = 8000; L = 9000; M = 8500; P = 142 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
